$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7480.5
$ws.Range("I28").Value = 8044.6924
$ws.Range("K28").Value = 8044.6924
$ws.Range("M28").Value = -7559.6924
$ws.Range("H31").Value = 35.75
$ws.Range("I31").Value = 35.75
$ws.Range("K31").Value = 107.25
$ws.Range("M31").Value = 122.75
$ws.Range("H48").Value = 1234.5
$ws.Range("J48").Value = 471
$ws.Range("L48").Value = 1413
$ws.Range("N48").Value = -1997
$ws.Range("H56").Value = 1234.5
$ws.Range("J56").Value = 471
$ws.Range("L56").Value = 1413
$ws.Range("N56").Value = -2481
$ws.Range("H62").Value = 5318.8
$ws.Range("I62").Value = 5099.4
$ws.Range("K62").Value = 5099.4
$ws.Range("M62").Value = -4475.4
$ws.Range("H65").Value = 5318.8
$ws.Range("I65").Value = 5099.4
$ws.Range("K65").Value = 25497
$ws.Range("M65").Value = -22377
$ws.Range("H132").Value = 3402.6
$ws.Range("I132").Value = 3217.7856
$ws.Range("K132").Value = 9653.356800000001
$ws.Range("M132").Value = -7123.356800000001
$ws.Range("H137").Value = 2247.5
$ws.Range("I137").Value = 2013.8572
$ws.Range("K137").Value = 6041.571599999999
$ws.Range("M137").Value = -3491.571599999999
$ws.Range("H141").Value = 2376.3845
$ws.Range("I141").Value = 2376.3845
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7129.1535
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1949.1535
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1439.8
$ws.Range("I25").Value = 549.75
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 549.75
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = -147.75
$ws.Range("N25").Value = -5804
$ws.Range("H30").Value = 600
$ws.Range("I30").Value = 600
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -450
$ws.Range("N30").Value = $null
$ws.Range("H35").Value = 3374.5
$ws.Range("I35").Value = 1999.5
$ws.Range("J35").Value = 4749.5
$ws.Range("K35").Value = 1999.5
$ws.Range("L35").Value = 4749.5
$ws.Range("M35").Value = -1593.5
$ws.Range("N35").Value = -5561.5
$ws.Range("H36").Value = 2805.5
$ws.Range("I36").Value = 2805.5
$ws.Range("K36").Value = 2805.5
$ws.Range("M36").Value = -2459.5
$ws.Range("H37").Value = 19285.715
$ws.Range("J37").Value = 23000
$ws.Range("L37").Value = 23000
$ws.Range("N37").Value = -23546
$ws.Range("H102").Value = 2158.5293
$ws.Range("I102").Value = 2003.2142
$ws.Range("K102").Value = 2003.2142
$ws.Range("M102").Value = -381.2141999999999
$ws.Range("H132").Value = 2424
$ws.Range("I132").Value = 2424
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7272
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4742
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("H94").Value = 3213.4285
$ws.Range("J94").Value = 2500
$ws.Range("L94").Value = 2500
$ws.Range("N94").Value = -3402
$ws.Range("H99").Value = 1278.75
$ws.Range("I99").Value = 934.5
$ws.Range("K99").Value = 934.5
$ws.Range("M99").Value = 563.5
$ws.Range("H105").Value = 1156.8
$ws.Range("I105").Value = 1124.2858
$ws.Range("J105").Value = 1232.6666
$ws.Range("K105").Value = 1124.2858
$ws.Range("L105").Value = 1232.6666
$ws.Range("M105").Value = 622.7141999999999
$ws.Range("N105").Value = -4726.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16540.334
$ws.Range("J28").Value = 16540.334
$ws.Range("L28").Value = 16540.334
$ws.Range("N28").Value = -17030.334
$ws.Range("H31").Value = 4701.1113
$ws.Range("I31").Value = 3901.5715
$ws.Range("J31").Value = 7499.5
$ws.Range("K31").Value = 3901.5715
$ws.Range("L31").Value = 7499.5
$ws.Range("M31").Value = -3606.5715
$ws.Range("N31").Value = -8089.5
$ws.Range("H34").Value = 4701.1113
$ws.Range("I34").Value = 3901.5715
$ws.Range("J34").Value = 7499.5
$ws.Range("K34").Value = 3901.5715
$ws.Range("L34").Value = 7499.5
$ws.Range("M34").Value = -3699.5715
$ws.Range("N34").Value = -7903.5
$ws.Range("H58").Value = 2465
$ws.Range("I58").Value = 2460.875
$ws.Range("K58").Value = 2460.875
$ws.Range("M58").Value = -2257.875
$ws.Range("H134").Value = 7280.2144
$ws.Range("I134").Value = 6926.9165
$ws.Range("K134").Value = 20780.7495
$ws.Range("M134").Value = -18245.7495
$ws.Range("H136").Value = 2465
$ws.Range("I136").Value = 2460.875
$ws.Range("K136").Value = 7382.625
$ws.Range("M136").Value = -4832.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 999.6667
$ws.Range("J25").Value = 999.5
$ws.Range("L25").Value = 2998.5
$ws.Range("N25").Value = -3336.5
$ws.Range("H29").Value = 360
$ws.Range("J29").Value = 400
$ws.Range("L29").Value = 1200
$ws.Range("N29").Value = -1754
$ws.Range("H30").Value = 999.6667
$ws.Range("J30").Value = 999.5
$ws.Range("L30").Value = 2998.5
$ws.Range("N30").Value = -3202.5
$ws.Range("H35").Value = 299
$ws.Range("I35").Value = 299
$ws.Range("K35").Value = 897
$ws.Range("M35").Value = -609
$ws.Range("H36").Value = 675
$ws.Range("J36").Value = 400
$ws.Range("L36").Value = 1200
$ws.Range("N36").Value = -1538
$ws.Range("H121").Value = 934.53845
$ws.Range("I121").Value = 348.5
$ws.Range("J121").Value = 1041.091
$ws.Range("K121").Value = 1045.5
$ws.Range("L121").Value = 3123.273
$ws.Range("M121").Value = 264.5
$ws.Range("N121").Value = -5743.272999999999
$ws.Range("H122").Value = 960.1429000000001
$ws.Range("J122").Value = 1245
$ws.Range("L122").Value = 11205
$ws.Range("N122").Value = -16105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 582.2273
$ws.Range("I2").Value = 601.75
$ws.Range("K2").Value = 601.75
$ws.Range("M2").Value = -488.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5500
$ws.Range("I4").Value = 5500
$ws.Range("K4").Value = 5500
$ws.Range("M4").Value = -5387
$ws.Range("H22").Value = 1799.8
$ws.Range("I22").Value = 1101
$ws.Range("J22").Value = 1974.5
$ws.Range("K22").Value = 1101
$ws.Range("L22").Value = 1974.5
$ws.Range("M22").Value = -806
$ws.Range("N22").Value = -2564.5
$ws.Range("H27").Value = 1799.8
$ws.Range("I27").Value = 1101
$ws.Range("J27").Value = 1974.5
$ws.Range("K27").Value = 1101
$ws.Range("L27").Value = 1974.5
$ws.Range("M27").Value = -994
$ws.Range("N27").Value = -2188.5
$ws.Range("H28").Value = 5500
$ws.Range("I28").Value = 5500
$ws.Range("K28").Value = 5500
$ws.Range("M28").Value = -5268
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("H37").Value = 5500
$ws.Range("I37").Value = 5500
$ws.Range("K37").Value = 5500
$ws.Range("M37").Value = -5393
$ws.Range("H55").Value = 141.375
$ws.Range("J55").Value = 59.2
$ws.Range("L55").Value = 59.2
$ws.Range("N55").Value = -405.2
$ws.Range("H68").Value = 499.66666
$ws.Range("I68").Value = 499.66666
$ws.Range("K68").Value = 499.66666
$ws.Range("M68").Value = 249.33334
$ws.Range("H71").Value = 499.66666
$ws.Range("I71").Value = 499.66666
$ws.Range("K71").Value = 2498.3333
$ws.Range("M71").Value = 1245.6667
$ws.Range("H74").Value = 90000
$ws.Range("I74").Value = 90000
$ws.Range("K74").Value = 90000
$ws.Range("M74").Value = -89002
$ws.Range("H77").Value = 90000
$ws.Range("I77").Value = 90000
$ws.Range("K77").Value = 270000
$ws.Range("M77").Value = -265008
$ws.Range("H80").Value = 30000
$ws.Range("I80").Value = 30000
$ws.Range("K80").Value = 30000
$ws.Range("M80").Value = -28877
$ws.Range("H83").Value = 30000
$ws.Range("I83").Value = 30000
$ws.Range("K83").Value = 90000
$ws.Range("M83").Value = -84384
$ws.Range("H132").Value = 8895.333000000001
$ws.Range("I132").Value = 9456.5
$ws.Range("K132").Value = 28369.5
$ws.Range("M132").Value = -25839.5
$ws.Range("H134").Value = 48999.2
$ws.Range("J134").Value = 48999.2
$ws.Range("L134").Value = 48999.2
$ws.Range("N134").Value = -59139.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 739.2
$ws.Range("I100").Value = 499.5
$ws.Range("J100").Value = 899
$ws.Range("K100").Value = 999
$ws.Range("L100").Value = 1798
$ws.Range("M100").Value = -458
$ws.Range("N100").Value = -2880
$ws.Range("H126").Value = 2988.4614
$ws.Range("I126").Value = 2685
$ws.Range("K126").Value = 8055
$ws.Range("M126").Value = -5585
